$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $range = $d.Content
    $ok = $range.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        Write-Output "WARNING: not found: $findText"
    }
    return $ok
}

# 1. R Core Team 2020 -> 2022 (intro paragraph citation)
Replace-Text "(R Core Team 2020)" "(R Core Team 2022)"

# 2. "An additional 4" -> "An additional 6"
Replace-Text "An additional 4 surveys" "An additional 6 surveys"

# 3. Insert the NCRO survey into the water-quality survey list
Replace-Text "DWR Stockton Dissolved Oxygen Survey (SDO), United States Bureau" "DWR Stockton Dissolved Oxygen Survey (SDO), DWR North Central Region Office monitoring (NCRO), United States Bureau"

# 4. Add NCRO to the nutrient variables survey list
Replace-Text "USGS_SFBS, USGS_CAWSC, and EMP surveys." "USGS_SFBS, USGS_CAWSC, EMP, and NCRO surveys."

# 5. CDEC) so -> CDEC), so  (also removes the gramStart/gramEnd proofErr pair around ")")
Replace-Text "monitoring program (CDEC) so these values" "monitoring program (CDEC), so these values"

# 6. before June 2019 so -> before June 2019, so (also removes the gramStart/gramEnd proofErr pair around "2019")
Replace-Text "before June 2019 so conductivity" "before June 2019, so conductivity"

# 7. "station,  latitude" (double space) -> "station, latitude" (single space); removes the gramStart/gramEnd proofErr pair
Replace-Text "date, time, station,  latitude, longitude" "date, time, station, latitude, longitude"

# 8. "(e.g. for non-fixed" -> "(e.g., for non-fixed"; removes the gramStart/gramEnd proofErr pair
Replace-Text "not available (e.g. for non-fixed" "not available (e.g., for non-fixed"

# 9. Reference list: R Core Team. 2020. -> R Core Team. 2022.
Replace-Text "R Core Team. 2020." "R Core Team. 2022."

Write-Output "All edits applied"
